# Add team record columns (Wins, Losses, Ties) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers, matching style of existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows (2 through 63) - same record (100-62-0) for every row.
$lastRow = 63
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 100  # AD
    $ws.Cells.Item($r, 31).Value = 62   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
